$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D2").Value = "26.933.38"
$ws.Range("E2").Value = "  +0.30%  "
$ws.Range("D3").Value = "1.554.38"
$ws.Range("E3").Value = "  +0.67%  "
$ws.Range("E4").Value = "  -0.04%  "
$ws.Range("D5").Value = "'206.62"
$ws.Range("E5").Value = "  +0.26%  "
$ws.Range("E6").Value = "  +0.45%  "
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("D8").Value = "'21.92"
$ws.Range("E8").Value = "  +2.34%  "
$ws.Range("E9").Value = "  -0.01%  "
$ws.Range("E10").Value = "  +0.67%  "
$ws.Range("D11").Value = "'0.0858"
$ws.Range("E11").Value = "  +0.44%  "
$ws.Range("D12").Value = "1.775.74"
$ws.Range("E12").Value = "  +0.72%  "
$ws.Range("D13").Value = "1.554.35"
$ws.Range("E13").Value = "  +0.67%  "
$ws.Range("E14").Value = "  +1.20%  "
$ws.Range("E15").Value = "  +1.27%  "
$ws.Range("D16").Value = "26.927.79"
$ws.Range("D17").Value = "'61.67"
$ws.Range("E17").Value = "  +0.60%  "
$ws.Range("D18").Value = "'217.92"
$ws.Range("E18").Value = "  +1.17%  "
$ws.Range("E19").Value = "  +1.61%  "
$ws.Range("E20").Value = "  +0.63%  "
$ws.Range("E21").Value = "  +0.00%  "
$ws.Range("E22").Value = "  +0.52%  "
$ws.Range("E23").Value = "  +0.11%  "
$ws.Range("D24").Value = "'1.94"
$ws.Range("E24").Value = "  +0.27%  "
$ws.Range("D25").Value = "'154.05"
$ws.Range("E25").Value = "  +0.67%  "
$ws.Range("E26").Value = "  -0.40%  "
$ws.Range("D27").Value = "'14.90"
$ws.Range("E27").Value = "  +0.27%  "
$ws.Range("E28").Value = "  +0.04%  "
$ws.Range("E29").Value = "  +0.56%  "
$ws.Range("E30").Value = "  +1.95%  "
$ws.Range("E31").Value = "  -0.35%  "
$ws.Range("E32").Value = "  -0.17%  "
$ws.Range("D33").Value = "1.434.69"
$ws.Range("E33").Value = "  +4.77%  "
$ws.Range("E34").Value = "  +4.06%  "
$ws.Range("E35").Value = "  +3.65%  "
$ws.Range("D36").Value = "'0.976"
$ws.Range("E36").Value = "  +1.56%  "
$ws.Range("E37").Value = "  +0.61%  "
$ws.Range("E38").Value = "  -0.11%  "
$ws.Range("D39").Value = "'0.520"
$ws.Range("E39").Value = "  -0.21%  "
$ws.Range("E40").Value = "  +0.27%  "
$ws.Range("E42").Value = "  -0.82%  "
$ws.Range("E43").Value = "  -0.42%  "
$ws.Range("E44").Value = "  +2.42%  "
$ws.Range("D45").Value = "'63.87"
$ws.Range("E45").Value = "  +0.95%  "
$ws.Range("E46").Value = "  +1.79%  "
$ws.Range("D47").Value = "1.689.52"
$ws.Range("E47").Value = "  +0.72%  "
$ws.Range("D48").Value = "'87.01"
$ws.Range("E48").Value = "  +3.35%  "
$ws.Range("E49").Value = "  +2.56%  "
$ws.Range("E50").Value = "  +3.36%  "
$ws.Range("D51").Value = "'0.0954"
$ws.Range("E51").Value = "  +1.46%  "
